$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Timp3"
$ws.Cells.Item(2,3).Value = "Agtr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 77.3174515
$ws.Cells.Item(2,8).Value = 154.634903
$ws.Cells.Item(2,9).Value = 0.3250124887881288
$ws.Cells.Item(2,10).Value = 0.2721424217793227
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.5
$ws.Cells.Item(2,13).Value = 0.1409635
$ws.Cells.Item(2,14).Value = 0.281927
$ws.Cells.Item(2,15).Value = 0.127302058387171
$ws.Cells.Item(2,16).Value = 0.0886289202516707
$ws.Cells.Item(2,17).Value = 10.89893857452025
$ws.Cells.Item(2,18).Value = 43.595754298081
$ws.Cells.Item(2,19).Value = 0.04137475882426613
$ws.Cells.Item(2,20).Value = 0.02411968899697613

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Timp3"
$ws.Cells.Item(3,3).Value = "Agtr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 77.3174515
$ws.Cells.Item(3,8).Value = 154.634903
$ws.Cells.Item(3,9).Value = 0.3250124887881288
$ws.Cells.Item(3,10).Value = 0.2721424217793227
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.9663516666666667
$ws.Cells.Item(3,14).Value = 2.899055
$ws.Cells.Item(3,15).Value = 0.872697941612829
$ws.Cells.Item(3,16).Value = 0.9113710797483293
$ws.Cells.Item(3,17).Value = 74.71584811944417
$ws.Cells.Item(3,18).Value = 448.2950887166651
$ws.Cells.Item(3,19).Value = 0.2836377299638627
$ws.Cells.Item(3,20).Value = 0.2480227327823466

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Timp3"
$ws.Cells.Item(4,3).Value = "Agtr2"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 13.28375266666667
$ws.Cells.Item(4,8).Value = 39.851258
$ws.Cells.Item(4,9).Value = 0.0558397286884097
$ws.Cells.Item(4,10).Value = 0.07013434646816191
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.5
$ws.Cells.Item(4,13).Value = 0.1409635
$ws.Cells.Item(4,14).Value = 0.281927
$ws.Cells.Item(4,15).Value = 0.127302058387171
$ws.Cells.Item(4,16).Value = 0.0886289202516707
$ws.Cells.Item(4,17).Value = 1.872524269027666
$ws.Cells.Item(4,18).Value = 11.235145614166
$ws.Cells.Item(4,19).Value = 0.007108512401815718
$ws.Cells.Item(4,20).Value = 0.006215931400029765

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Timp3"
$ws.Cells.Item(5,3).Value = "Agtr2"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 13.28375266666667
$ws.Cells.Item(5,8).Value = 39.851258
$ws.Cells.Item(5,9).Value = 0.0558397286884097
$ws.Cells.Item(5,10).Value = 0.07013434646816191
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.9663516666666667
$ws.Cells.Item(5,14).Value = 2.899055
$ws.Cells.Item(5,15).Value = 0.872697941612829
$ws.Cells.Item(5,16).Value = 0.9113710797483293
$ws.Cells.Item(5,17).Value = 12.83677652902111
$ws.Cells.Item(5,18).Value = 115.53098876119
$ws.Cells.Item(5,19).Value = 0.04873121628659398
$ws.Cells.Item(5,20).Value = 0.06391841506813215

# Row 6
$ws.Cells.Item(6,1).Value = "Neutro"
$ws.Cells.Item(6,2).Value = "Timp3"
$ws.Cells.Item(6,3).Value = "Agtr2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 79.14797966666667
$ws.Cells.Item(6,8).Value = 237.443939
$ws.Cells.Item(6,9).Value = 0.332707317105706
$ws.Cells.Item(6,10).Value = 0.4178782884241973
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.5
$ws.Cells.Item(6,13).Value = 0.1409635
$ws.Cells.Item(6,14).Value = 0.281927
$ws.Cells.Item(6,15).Value = 0.127302058387171
$ws.Cells.Item(6,16).Value = 0.0886289202516707
$ws.Cells.Item(6,17).Value = 11.15697623174217
$ws.Cells.Item(6,18).Value = 66.941857390453
$ws.Cells.Item(6,19).Value = 0.0423543263080296
$ws.Cells.Item(6,20).Value = 0.03703610149965283

# Row 7
$ws.Cells.Item(7,1).Value = "Neutro"
$ws.Cells.Item(7,2).Value = "Timp3"
$ws.Cells.Item(7,3).Value = "Agtr2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 79.14797966666667
$ws.Cells.Item(7,8).Value = 237.443939
$ws.Cells.Item(7,9).Value = 0.332707317105706
$ws.Cells.Item(7,10).Value = 0.4178782884241973
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.9663516666666667
$ws.Cells.Item(7,14).Value = 2.899055
$ws.Cells.Item(7,15).Value = 0.872697941612829
$ws.Cells.Item(7,16).Value = 0.9113710797483293
$ws.Cells.Item(7,17).Value = 76.48478206418278
$ws.Cells.Item(7,18).Value = 688.363038577645
$ws.Cells.Item(7,19).Value = 0.2903529907976764
$ws.Cells.Item(7,20).Value = 0.3808421869245445

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Timp3"
$ws.Cells.Item(8,3).Value = "Agtr2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 68.1415255
$ws.Cells.Item(8,8).Value = 136.283051
$ws.Cells.Item(8,9).Value = 0.2864404654177555
$ws.Cells.Item(8,10).Value = 0.239844943328318
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.5
$ws.Cells.Item(8,13).Value = 0.1409635
$ws.Cells.Item(8,14).Value = 0.281927
$ws.Cells.Item(8,15).Value = 0.127302058387171
$ws.Cells.Item(8,16).Value = 0.0886289202516707
$ws.Cells.Item(8,17).Value = 9.60546792981925
$ws.Cells.Item(8,18).Value = 38.421871719277
$ws.Cells.Item(8,19).Value = 0.03646446085305954
$ws.Cells.Item(8,20).Value = 0.02125719835501197

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Timp3"
$ws.Cells.Item(9,3).Value = "Agtr2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 68.1415255
$ws.Cells.Item(9,8).Value = 136.283051
$ws.Cells.Item(9,9).Value = 0.2864404654177555
$ws.Cells.Item(9,10).Value = 0.239844943328318
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.9663516666666667
$ws.Cells.Item(9,14).Value = 2.899055
$ws.Cells.Item(9,15).Value = 0.872697941612829
$ws.Cells.Item(9,16).Value = 0.9113710797483293
$ws.Cells.Item(9,17).Value = 65.84867673613417
$ws.Cells.Item(9,18).Value = 395.092060416805
$ws.Cells.Item(9,19).Value = 0.249976004564696
$ws.Cells.Item(9,20).Value = 0.218587744973306

